# Auto-generated script to update cryptos.xlsx per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "53.240.61"
$ws.Range("E2").Value = "  -4.98%  "

$ws.Range("D3").Value = "2.196.06"
$ws.Range("E3").Value = "  -7.12%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "480.06"
$ws.Range("E5").Value = "  -4.11%  "

$ws.Range("D6").Value = "124.44"
$ws.Range("E6").Value = "  -3.08%  "

$ws.Range("E7").Value = "  +0.23%  "

$ws.Range("D8").Value = "0.516"
$ws.Range("E8").Value = "  -5.09%  "

$ws.Range("D9").Value = "2.196.33"
$ws.Range("E9").Value = "  -7.21%  "

$ws.Range("E10").Value = "  -1.63%  "

$ws.Range("D11").Value = "0.0905"
$ws.Range("E11").Value = "  -7.52%  "

$ws.Range("D12").Value = "4.64"
$ws.Range("E12").Value = "  -2.91%  "

$ws.Range("D13").Value = "0.311"
$ws.Range("E13").Value = "  -3.39%  "

$ws.Range("D14").Value = "2.597.82"
$ws.Range("E14").Value = "  -6.75%  "

$ws.Range("D15").Value = "20.84"
$ws.Range("E15").Value = "  -2.56%  "

$ws.Range("D16").Value = "53.258.15"
$ws.Range("E16").Value = "  -4.90%  "

$ws.Range("D17").Value = "0.0000126"
$ws.Range("E17").Value = "  -4.18%  "

$ws.Range("D18").Value = "2.211.08"
$ws.Range("E18").Value = "  -7.10%  "

$ws.Range("D19").Value = "9.45"
$ws.Range("E19").Value = "  -5.29%  "

$ws.Range("D20").Value = "3.91"
$ws.Range("E20").Value = "  -3.03%  "

$ws.Range("D21").Value = "294.85"
$ws.Range("E21").Value = "  -3.70%  "

$ws.Range("D22").Value = "6.02"
$ws.Range("E22").Value = "  -4.04%  "

$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D24").Value = "62.75"
$ws.Range("E24").Value = "  -4.73%  "

$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.11%  "

$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").Value = "0.360"
$ws.Range("E26").Value = "  -1.77%  "

$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.315.32"
$ws.Range("E27").Value = "  -6.42%  "

$ws.Range("E28").Value = "  -3.09%  "

$ws.Range("D29").Value = "6.91"
$ws.Range("E29").Value = "  -3.74%  "

$ws.Range("D30").Value = "167.38"
$ws.Range("E30").Value = "  -2.80%  "

$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("D32").Value = "1.56"
$ws.Range("E32").Value = "  -4.32%  "

$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "0.997"
$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").Value = "0.0₃0663"
$ws.Range("E34").Value = "  -6.73%  "

$ws.Range("D35").Value = "5.62"
$ws.Range("E35").Value = "  -1.94%  "

$ws.Range("E36").Value = "  -3.68%  "

$ws.Range("D37").Value = "17.16"
$ws.Range("E37").Value = "  -2.48%  "

$ws.Range("D38").Value = "1.13"
$ws.Range("E38").Value = "  -3.47%  "

$ws.Range("D39").Value = "0.819"
$ws.Range("E39").Value = "  +4.18%  "

$ws.Range("D40").Value = "3.52"
$ws.Range("E40").Value = "  -5.67%  "

$ws.Range("D41").Value = "35.63"
$ws.Range("E41").Value = "  -1.18%  "

$ws.Range("D42").Value = "0.362"
$ws.Range("E42").Value = "  -1.63%  "

$ws.Range("D43").Value = "1.34"
$ws.Range("E43").Value = "  -2.36%  "

$ws.Range("D44").Value = "3.22"
$ws.Range("E44").Value = "  -3.83%  "

$ws.Range("D45").Value = "121.12"
$ws.Range("E45").Value = "  -6.50%  "

$ws.Range("D46").Value = "4.59"
$ws.Range("E46").Value = "  -1.57%  "

$ws.Range("D47").Value = "0.0870"
$ws.Range("E47").Value = "  -3.28%  "

$ws.Range("D48").Value = "0.526"
$ws.Range("E48").Value = "  -6.38%  "

$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "0.0466"
$ws.Range("E49").Value = "  -2.71%  "

$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").Value = "227.24"
$ws.Range("E50").Value = "  -5.22%  "

$ws.Range("D51").Value = "0.0199"
$ws.Range("E51").Value = "  -3.43%  "
